# Leave Card update — insert two new travel/leave rows (3/8-12/2024 ITALY,
# 3/27-31/2024 THAILAND) into Table1 on Sheet1 right after row 45, which
# pushes every following data row down by two and grows the table by two
# rows at the bottom (new blank rows keep the table's recurring
# month-period date sequence intact).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("CONVERTION")
$lo  = $ws1.ListObjects.Item(1)

$formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# 1. Insert two blank worksheet rows right before the current row 46 — this
#    shifts rows 46..144 down to 48..146 and keeps every formula/value in
#    those rows intact (dates, formulas, remarks, etc. all move as a unit).
$ws1.Range("A46:A47").EntireRow.Insert()

# 2. Grow the table definition to match (A8:K144 -> A8:K146).
$lo.Resize($ws1.Range("A8:K146"))

# 3. The two freshly inserted rows have generic (non-table) styling; clone
#    the look of a normal data row (row 48, the old row 46) into them so
#    borders/number-formats match the rest of Table1.
$ws1.Range("A48:K48").Copy($ws1.Range("A46:K46"))
$ws1.Range("A48:K48").Copy($ws1.Range("A47:K47"))
$excel.CutCopyMode = $false

# Re-assert the calculated "EARNED " column formula on every row the Copy
# above touched (Copy in this host does not carry the formula across).
$ws1.Range("G46").Formula = $formula
$ws1.Range("G47").Formula = $formula

# 4. Populate the two new rows with the new leave records. Row 46 = VL(3-0-0)
#    trip (3 days, "3/8-12/2024 ITALY"); row 47 = VL(1-0-0) trip (1 day,
#    "3/27-31/2024 THAILAND"). Column A (date-of-period) is left blank on
#    these two rows, same as the other in-between rows.
$ws1.Range("A46").ClearContents()
$ws1.Range("B46").Value = "VL(3-0-0)"
$ws1.Range("F46").Value = 3
$ws1.Range("K46").Value = "3/8-12/2024 ITALY"

$ws1.Range("A47").ClearContents()
$ws1.Range("B47").Value = "VL(1-0-0)"
$ws1.Range("F47").Value = 1
$ws1.Range("K47").Value = "3/27-31/2024 THAILAND"

# 5. The table's calculated "EARNED " column formula on the two brand new
#    rows at the bottom of the grown table (145 & 146) needs to be the full
#    structured-reference form too (matches every other row in the column).
$ws1.Range("G145").Formula = $formula
$ws1.Range("G146").Formula = $formula

# 6. Leave the cursor where the author left it (row 46 was the former
#    selection; after the insert the same logical spot is row 48).
$ws1.Range("F48").Select()

# 7. Scroll the CONVERTION sheet down a little (cosmetic) while keeping its
#    existing selection and the workbook's active sheet/tab unchanged.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 20
$ws2.Range("J34").Select()
$ws1.Activate()
